$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value2 = 88
$ws.Range("F5").Value2 = 1042
$ws.Range("F7").Value2 = 2684
$ws.Range("F8").Value2 = 40
$ws.Range("F9").Value2 = 1313
$ws.Range("F11").Value2 = 634
$ws.Range("F12").Value2 = 947
$ws.Range("F13").Value2 = 1194
$ws.Range("F14").Value2 = 298
$ws.Range("F16").Value2 = 751
$ws.Range("F17").Value2 = 798
$ws.Range("F18").Value2 = 226
$ws.Range("F19").Value2 = 538
$ws.Range("F20").Value2 = 1145
$ws.Range("F22").Value2 = 654
$ws.Range("F23").Value2 = 614
$ws.Range("F24").Value2 = 234
$ws.Range("F25").Value2 = 320
$ws.Range("F27").Value2 = 700
$ws.Range("F28").Value2 = 604
$ws.Range("F29").Value2 = 5918
$ws.Range("F30").Value2 = 503
$ws.Range("F34").Value2 = 186
$ws.Range("F37").Value2 = 108
$ws.Range("F38").Value2 = 452
$ws.Range("F39").Value2 = 147
$ws.Range("F41").Value2 = 155
$ws.Range("F47").Value2 = 123

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value2 = 101
$ws.Range("F9").Value2 = 54
$ws.Range("F12").Value2 = 198
$ws.Range("F13").Value2 = 4414
$ws.Range("F14").Value2 = 41
$ws.Range("F17").Value2 = 41
$ws.Range("F20").Value2 = 1

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 753

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 753
$ws.Range("F6").Value2 = 1042
$ws.Range("F7").Value2 = 2685
$ws.Range("F8").Value2 = 40
$ws.Range("F9").Value2 = 1313
$ws.Range("F11").Value2 = 634
$ws.Range("F12").Value2 = 947
$ws.Range("F13").Value2 = 1194
$ws.Range("F14").Value2 = 298
$ws.Range("F17").Value2 = 751
$ws.Range("F19").Value2 = 798
$ws.Range("F20").Value2 = 226
$ws.Range("F21").Value2 = 538
$ws.Range("F22").Value2 = 1145
$ws.Range("F24").Value2 = 54
$ws.Range("F25").Value2 = 654
$ws.Range("F26").Value2 = 614
$ws.Range("F27").Value2 = 234
$ws.Range("F28").Value2 = 320
$ws.Range("F30").Value2 = 604
$ws.Range("F31").Value2 = 5922
$ws.Range("F32").Value2 = 198
$ws.Range("F33").Value2 = 503
$ws.Range("F36").Value2 = 186
$ws.Range("F39").Value2 = 452
$ws.Range("F40").Value2 = 41
$ws.Range("F41").Value2 = 41
$ws.Range("F43").Value2 = 41
$ws.Range("F48").Value2 = 123
